# Generate Report for Handback
# Update the "generate date" / handoff / handback timestamps that are
# refreshed each time the handback report is regenerated.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# zh-cn sheet (row 4 = f87b9357-...) : Correspond Handoff / Handback Datetime
$zhcn.Range("H4").Value = "2016-08-26 22:44:13"
$zhcn.Range("K4").Value = "2016-08-26 22:44:31"

# de-de sheet (row 4 = f87b9357-...) : Correspond Handoff / Handback Datetime
$dede.Range("H4").Value = "2016-08-26 22:44:17"
$dede.Range("K4").Value = "2016-08-26 22:44:37"

# Overview sheet (row 4 = f87b9357-...) : Latest HO Xliff Generate Date
# (mirrors the de-de handoff datetime above)
$overview.Range("G4").Value = "2016-08-26 22:44:17"
